# Add a "time_taken" column (F) to the worksheet, recording a processing
# timestamp for each data row - panelApp metadata update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: text + formatting copied from the adjacent header cell (E1)
# so it matches the bold/bordered/centered header style (s="1").
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

$timestamps = @(
    "2021-10-05 13:40:47.770019",
    "2021-10-05 13:40:47.770030",
    "2021-10-05 13:40:47.770034",
    "2021-10-05 13:40:47.770036",
    "2021-10-05 13:40:47.770039",
    "2021-10-05 13:40:47.770042",
    "2021-10-05 13:40:47.770045",
    "2021-10-05 13:40:47.770047",
    "2021-10-05 13:40:47.770050",
    "2021-10-05 13:40:47.770052",
    "2021-10-05 13:40:47.770055",
    "2021-10-05 13:40:47.770057",
    "2021-10-05 13:40:47.770060",
    "2021-10-05 13:40:47.770062",
    "2021-10-05 13:40:47.770065",
    "2021-10-05 13:40:47.770067",
    "2021-10-05 13:40:47.770070",
    "2021-10-05 13:40:47.770073",
    "2021-10-05 13:40:47.770075",
    "2021-10-05 13:40:47.770078",
    "2021-10-05 13:40:47.770080",
    "2021-10-05 13:40:47.770083",
    "2021-10-05 13:40:47.770085",
    "2021-10-05 13:40:47.770088",
    "2021-10-05 13:40:47.770091",
    "2021-10-05 13:40:47.770094",
    "2021-10-05 13:40:47.770096",
    "2021-10-05 13:40:47.770099",
    "2021-10-05 13:40:47.770101",
    "2021-10-05 13:40:47.770104",
    "2021-10-05 13:40:47.770106",
    "2021-10-05 13:40:47.770109",
    "2021-10-05 13:40:47.770112",
    "2021-10-05 13:40:47.770115",
    "2021-10-05 13:40:47.770117",
    "2021-10-05 13:40:47.770120",
    "2021-10-05 13:40:47.770122",
    "2021-10-05 13:40:47.770125",
    "2021-10-05 13:40:47.770128",
    "2021-10-05 13:40:47.770130",
    "2021-10-05 13:40:47.770133",
    "2021-10-05 13:40:47.770136",
    "2021-10-05 13:40:47.770138",
    "2021-10-05 13:40:47.770141",
    "2021-10-05 13:40:47.770143"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
